$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 (buyRate): "null" -> "0.93"
$ws.Range("D2").Formula = '="0.93"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

# E2 (midRate): "0.8496" -> "null"
$ws.Range("E2").Value = "null"

# H2 (rateEffectiveDate): "2019-12-16" -> "2020-11-25"
$ws.Range("H2").Formula = '="2020-11-25"'
$ws.Range("H2").Copy()
$ws.Range("H2").PasteSpecial(-4163)

# K2 (rateTenor): "001M" -> "003M"
$ws.Range("K2").Value = "003M"

$excel.CutCopyMode = $false
